$wb = $excel.ActiveWorkbook

# Sheet: Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 6004
$ws.Range("J3").Value = 6410
$ws.Range("H4").Value = 1703
$ws.Range("J4").Value = 1386
$ws.Range("J5").Value = 491
$ws.Range("J6").Value = 8293
$ws.Range("H7").Value = 26014
$ws.Range("J7").Value = 22584

# Sheet: By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J6").Value = 166
$ws.Range("J7").Value = 669
$ws.Range("J8").Value = 1420
$ws.Range("J9").Value = 115
$ws.Range("J11").Value = 360
$ws.Range("J14").Value = 113
$ws.Range("J15").Value = 250
$ws.Range("J19").Value = 668
$ws.Range("J20").Value = 473
$ws.Range("J23").Value = 212
$ws.Range("J27").Value = 139
$ws.Range("J29").Value = 1243
$ws.Range("J33").Value = 1042
$ws.Range("J37").Value = 691
$ws.Range("J42").Value = 954
$ws.Range("J43").Value = 187
$ws.Range("J46").Value = 74
$ws.Range("J48").Value = 265
$ws.Range("J50").Value = 134
$ws.Range("J53").Value = 314
$ws.Range("J55").Value = 318
$ws.Range("J57").Value = 101
$ws.Range("J60").Value = 133
$ws.Range("H63").Value = 261
$ws.Range("I63").Value = 247
$ws.Range("J63").Value = 81
$ws.Range("J64").Value = 149
$ws.Range("J65").Value = 561
$ws.Range("J67").Value = 858
$ws.Range("J68").Value = 45
$ws.Range("J76").Value = 344
$ws.Range("J77").Value = 169
$ws.Range("J79").Value = 642
$ws.Range("I86").Value = 168
$ws.Range("J87").Value = 76
$ws.Range("J88").Value = 237
$ws.Range("J89").Value = 295
$ws.Range("J91").Value = 258
$ws.Range("J93").Value = 99
$ws.Range("J94").Value = 232
$ws.Range("J96").Value = 257
$ws.Range("H101").Value = 26014
$ws.Range("J101").Value = 22584

# Sheet: Bridgeport
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J4").Value = 7
$ws.Range("J7").Value = 113

# Sheet: West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 77
$ws.Range("J7").Value = 257

# Sheet: Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 205
$ws.Range("J7").Value = 669

# Sheet: Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 108
$ws.Range("J6").Value = 154
$ws.Range("J7").Value = 360

# Sheet: Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J4").Value = 29
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 295

# Sheet: Logan Square
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 56
$ws.Range("J3").Value = 40
$ws.Range("J6").Value = 207
$ws.Range("J7").Value = 314

# Sheet: Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J3").Value = 433
$ws.Range("J6").Value = 489
$ws.Range("J7").Value = 1420

# Sheet: Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 246
$ws.Range("J3").Value = 345
$ws.Range("J6").Value = 364
$ws.Range("J7").Value = 1042

# Sheet: Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 205
$ws.Range("J3").Value = 237
$ws.Range("J6").Value = 199
$ws.Range("J7").Value = 691

# Sheet: New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 163
$ws.Range("J7").Value = 561

# Sheet: North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 324
$ws.Range("J6").Value = 231
$ws.Range("J7").Value = 858

# Sheet: Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 383
$ws.Range("J3").Value = 432
$ws.Range("J6").Value = 314
$ws.Range("J7").Value = 1243

# Sheet: Lake View
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J4").Value = 40
$ws.Range("J7").Value = 265

# Sheet: Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J6").Value = 258
$ws.Range("J7").Value = 668

# Sheet: River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 191
$ws.Range("J7").Value = 344

# Sheet: Ashburn
$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J2").Value = 49
$ws.Range("J7").Value = 166

# Sheet: Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 204
$ws.Range("J7").Value = 954

# Sheet: Lower West Side
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J3").Value = 68
$ws.Range("J7").Value = 318

# Sheet: Jefferson Park
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 74

# Sheet: Douglas
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 212

# Sheet: Washington Park
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 72
$ws.Range("J7").Value = 258

# Sheet: Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J6").Value = 187
$ws.Range("J7").Value = 642

# Sheet: Near South Side
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J3").Value = 39
$ws.Range("J7").Value = 149

# Sheet: Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J5").Value = 11
$ws.Range("J7").Value = 473

# Sheet: West Lawn
$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J2").Value = 24
$ws.Range("J7").Value = 99

# Sheet: West Loop
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 129
$ws.Range("J7").Value = 232

# Sheet: Brighton Park
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J2").Value = 74
$ws.Range("J7").Value = 250

# Sheet: Lincoln Square
$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J4").Value = 21
$ws.Range("J7").Value = 134

# Sheet: Avalon Park
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J2").Value = 32
$ws.Range("J7").Value = 115

# Sheet: United Center
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J3").Value = 63
$ws.Range("J7").Value = 237

# Sheet: Edgewater
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 32
$ws.Range("J4").Value = 18
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 139

# Sheet: Streeterville
$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 79
$ws.Range("I7").Value = 168

# Sheet: North Park
$ws = $wb.Worksheets.Item("North Park")
$ws.Range("J2").Value = 18
$ws.Range("J7").Value = 45

# Sheet: Mckinley Park
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J3").Value = 26
$ws.Range("J7").Value = 101

# Sheet: Morgan Park
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J2").Value = 48
$ws.Range("J7").Value = 133

# Sheet: Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J2").Value = 21
$ws.Range("J6").Value = 110
$ws.Range("J7").Value = 187

# Sheet: Riverdale
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 169

# Sheet: Ukrainian Village
$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 76
